$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal-string value into a cell, forcing text storage
# (leading apostrophe keeps values like "True"/"False"/"" from being
# coerced into booleans, and keeps everything as shared-string text).
# ---------------------------------------------------------------------------
function Set-TextCell($cell, [string]$text) {
    if ($text -eq "") {
        $cell.Value() = "'"
    } else {
        $cell.Value() = "'" + $text
    }
}

# ===========================================================================
# Sheet "Overview" (sheet1): two new rows for the handed-back files
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

$overviewRows = @(
    @{ Row = 6; Name = "d6f54bbc-b386-4234-8df3-09507feba5f4.md"; HoDate = "2016-09-01 06:46:51" },
    @{ Row = 7; Name = "f895d7e5-fbe1-47f9-b01c-59ff9d330680.md"; HoDate = "2016-09-01 06:46:51" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $name = $r.Name
    $display = "e2e\" + $name

    Set-TextCell $wsOverview.Cells.Item($row, 1) $name
    Set-TextCell $wsOverview.Cells.Item($row, 3) ".md"
    Set-TextCell $wsOverview.Cells.Item($row, 4) ""
    Set-TextCell $wsOverview.Cells.Item($row, 5) "Ready for handoff"
    Set-TextCell $wsOverview.Cells.Item($row, 6) "Ready for handoff"
    Set-TextCell $wsOverview.Cells.Item($row, 7) $r.HoDate

    $wsOverview.Cells.Item($row, 2).Value() = "'" + $display
    $linkCell = $wsOverview.Cells.Item($row, 2)
    $wsOverview.Hyperlinks.Add($linkCell, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/$name", "", "", $display) | Out-Null
    $linkCell.Font.Underline = 2
    $linkCell.Font.Color = 15570276
}

# ===========================================================================
# Sheets "zh-cn" (sheet2) and "de-de" (sheet3): two new rows each
# ===========================================================================
$langSheets = @(
    @{
        Name = "zh-cn";
        Files = @(
            @{ Row = 6; Name = "d6f54bbc-b386-4234-8df3-09507feba5f4.md"; Xlf = "d6f54bbc-b386-4234-8df3-09507feba5f4.8b3a43ec701b81a713e992bcee41ff1ddde1b7f0.zh-cn.xlf"; HoDate = "2016-09-01 06:46:47" },
            @{ Row = 7; Name = "f895d7e5-fbe1-47f9-b01c-59ff9d330680.md"; Xlf = "f895d7e5-fbe1-47f9-b01c-59ff9d330680.d9f0e7268bcfb2c6e6669fe1da41d38824c7e65e.zh-cn.xlf"; HoDate = "2016-09-01 06:46:47" }
        )
    },
    @{
        Name = "de-de";
        Files = @(
            @{ Row = 6; Name = "d6f54bbc-b386-4234-8df3-09507feba5f4.md"; Xlf = "d6f54bbc-b386-4234-8df3-09507feba5f4.8b3a43ec701b81a713e992bcee41ff1ddde1b7f0.de-de.xlf"; HoDate = "2016-09-01 06:46:51" },
            @{ Row = 7; Name = "f895d7e5-fbe1-47f9-b01c-59ff9d330680.md"; Xlf = "f895d7e5-fbe1-47f9-b01c-59ff9d330680.d9f0e7268bcfb2c6e6669fe1da41d38824c7e65e.de-de.xlf"; HoDate = "2016-09-01 06:46:51" }
        )
    }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)
    $lo = $ws.ListObjects.Item(1)
    $lo.ListRows.Add() | Out-Null
    $lo.ListRows.Add() | Out-Null

    foreach ($f in $lang.Files) {
        $row = $f.Row
        $name = $f.Name

        $wsCellA = $ws.Cells.Item($row, 1)
        $wsCellA.Value() = "'" + $name
        $ws.Hyperlinks.Add($wsCellA, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/$name", "", "", $name) | Out-Null
        $wsCellA.Font.Underline = 2
        $wsCellA.Font.Color = 15570276

        Set-TextCell $ws.Cells.Item($row, 2) ".md"
        Set-TextCell $ws.Cells.Item($row, 3) "Ready for handoff"
        Set-TextCell $ws.Cells.Item($row, 4) "e2e"
        Set-TextCell $ws.Cells.Item($row, 5) "ht"
        Set-TextCell $ws.Cells.Item($row, 6) "False"
        Set-TextCell $ws.Cells.Item($row, 7) $f.Xlf
        Set-TextCell $ws.Cells.Item($row, 8) $f.HoDate
        Set-TextCell $ws.Cells.Item($row, 9) ""
        Set-TextCell $ws.Cells.Item($row, 10) ""
        Set-TextCell $ws.Cells.Item($row, 11) "0001-01-01 00:00:00"
        Set-TextCell $ws.Cells.Item($row, 12) ""
        Set-TextCell $ws.Cells.Item($row, 13) "True"
        Set-TextCell $ws.Cells.Item($row, 14) ""
        Set-TextCell $ws.Cells.Item($row, 15) "False"
        Set-TextCell $ws.Cells.Item($row, 16) ""
    }
}
